$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column C (shifts old C..K to D..L)
$ws.Columns.Item(3).Insert()

# Update header row
$ws.Range("B1").Value = "Denomination"
$ws.Range("C1").Value = "taille"

# Fill in the new "taille" (size/denomination) column
$ws.Range("C2").Value = "6"
$ws.Range("C3").Value = "8"
$ws.Range("C4").Value = "10"
$ws.Range("C5").Value = "12"
$ws.Range("C6").Value = "15"
$ws.Range("C7").Value = "18"
$ws.Range("C8").Value = "22"
$ws.Range("C9").Value = "28"
$ws.Range("C10").Value = "12-10"
$ws.Range("C11").Value = "15-10"
$ws.Range("C12").Value = "15-12"
$ws.Range("C13").Value = "6"
$ws.Range("C14").Value = "8"
$ws.Range("C15").Value = "10"
$ws.Range("C16").Value = "12"
$ws.Range("C17").Value = "15"
$ws.Range("C18").Value = "18"
$ws.Range("C19").Value = "22"
$ws.Range("C20").Value = "28"
$ws.Range("C21").Value = "6"
$ws.Range("C22").Value = "8"
$ws.Range("C23").Value = "10"
$ws.Range("C24").Value = "12"
$ws.Range("C25").Value = "15"
$ws.Range("C26").Value = "18"
$ws.Range("C27").Value = "22"
$ws.Range("C28").Value = "28"
$ws.Range("C29").Value = "1/2F-12"
$ws.Range("C30").Value = "3/4F-15"
$ws.Range("C31").Value = "3/8F-10"
$ws.Range("C32").Value = "3/8F-12"
$ws.Range("C33").Value = "1/2F-15"
$ws.Range("C34").Value = "3/4F-22"
$ws.Range("C35").Value = "3/4F-18"
$ws.Range("C36").Value = "1/2F-10"
$ws.Range("C37").Value = "1/2M-12"
$ws.Range("C38").Value = "1/2M-15"
$ws.Range("C39").Value = "3/4M-15"
$ws.Range("C40").Value = "3/4M-18"
$ws.Range("C41").Value = "3/4M-22"
$ws.Range("C42").Value = "6"
$ws.Range("C43").Value = "8"
$ws.Range("C44").Value = "6"
$ws.Range("C45").Value = "8"
$ws.Range("C46").Value = "10"
$ws.Range("C47").Value = "12"
$ws.Range("C48").Value = "15"
$ws.Range("C49").Value = "10"
$ws.Range("C50").Value = "12"
$ws.Range("C51").Value = "15"
$ws.Range("C52").Value = "15"
$ws.Range("C53").Value = "18"
$ws.Range("C54").Value = "22"
$ws.Range("C55").Value = "22"
$ws.Range("C56").Value = "28"
$ws.Range("C57").Value = "6"
$ws.Range("C58").Value = "8"
$ws.Range("C59").Value = "6"
$ws.Range("C60").Value = "8"
$ws.Range("C61").Value = "10"
$ws.Range("C62").Value = "12"
$ws.Range("C63").Value = "15"
$ws.Range("C64").Value = "10"
$ws.Range("C65").Value = "12"
$ws.Range("C66").Value = "15"
$ws.Range("C67").Value = "15"
$ws.Range("C68").Value = "18"
$ws.Range("C69").Value = "22"
$ws.Range("C70").Value = "28"
$ws.Range("C71").Value = "1/-15"
$ws.Range("C72").Value = "1/-12"

# Column C width (narrow, best-fit like the other text columns)
$ws.Columns.Item(3).ColumnWidth = 7

# Restore selection/scroll position recorded by Excel on save
$ws.Range("C31").Select() | Out-Null
